$d = $word.ActiveDocument

function Strike-Paragraph([string]$oldText) {
    # Step 1: force adjacent runs with identical resulting formatting to merge
    # by replacing the text with itself via Find/Replace.
    $merge = $d.Content.Duplicate
    $merge.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2) | Out-Null

    # Step 2: find the (now merged) run and apply strikethrough formatting.
    $hit = $d.Content.Duplicate
    $hit.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
    $hit.Font.StrikeThrough = 1

    # Step 3: also strike the paragraph mark itself, matching Word's behaviour
    # when the whole paragraph is selected and struck through.
    $mark = $hit.Duplicate
    $mark.Collapse(0)
    $mark.Font.StrikeThrough = 1
}

# 1) New "Key:" paragraph at the very top of the document.
$firstPara = $d.Paragraphs(1).Range
$firstPara.InsertParagraphBefore() | Out-Null
$keyPara = $d.Paragraphs(1)
$keyPara.Range.Text = "Key: "
$keyPara.Range.InsertAfter("line through") | Out-Null
$keyPara.Range.InsertAfter(" = Done") | Out-Null

$keyHit = $d.Content.Duplicate
$keyHit.Find.Execute("line through", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$keyHit.Font.StrikeThrough = 1

# 2) Ability ball mechanics bullets that are now "done".
Strike-Paragraph("A ball that can be hit without any special mechanics (similarly to how the player hits the target ball in the initial stage of the game)")
Strike-Paragraph("Ball that sticks to the surface on first collision and allows the player to aim the angle of reflection.")

# 3) Phase-through-wall bullet: strike existing text, then append the
#    "(Altered: ...)" note as its own run, also struck through.
$phaseOld = "A ball which can phase through the first wall it collides with and bounces off every wall after that."
$phaseNote = "(Altered: can go through a specific wall block, and as many as it collides with)"

$phaseMerge = $d.Content.Duplicate
$phaseMerge.Find.Execute($phaseOld, $true, $false, $false, $false, $false, $true, 1, $false, $phaseOld, 2) | Out-Null

$phaseHit = $d.Content.Duplicate
$phaseHit.Find.Execute($phaseOld, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$phaseHit.Font.StrikeThrough = 1
$phaseHit.Collapse(0)
$phaseHit.InsertAfter($phaseNote) | Out-Null

$phaseNoteHit = $d.Content.Duplicate
$phaseNoteHit.Find.Execute($phaseNote, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$phaseNoteHit.Font.StrikeThrough = 1

$phaseMark = $phaseNoteHit.Duplicate
$phaseMark.Collapse(0)
$phaseMark.Font.StrikeThrough = 1

# 4) Level object bullets that are now "done".
Strike-Paragraph("Static wall")
Strike-Paragraph("Dead zone or fall void that ball falls off of")
Strike-Paragraph("Angled wall or triangle")
Strike-Paragraph("Ramp")

# 5) UI interface bullets that are now "done".
Strike-Paragraph("Current ability of power ball")
Strike-Paragraph("Available selection of power ball abilities")
Strike-Paragraph("Shot counter")
Strike-Paragraph("Par of level")
